$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 157. This shifts existing rows 157-278 down to 158-279,
# carrying their formatting (and thus the date style on column D) along with them.
$ws.Rows.Item(157).Insert()

# Populate the newly inserted row 157 with the new record.
$ws.Cells.Item(157, 1).Value = 4
$ws.Cells.Item(157, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(157, 3).Value = "Los Lagos"
$ws.Cells.Item(157, 4).Value = 44651
$ws.Cells.Item(157, 5).Value = 10
$ws.Cells.Item(157, 6).Value = 100114014
$ws.Cells.Item(157, 7).Value = "Betarraga"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 500
$ws.Cells.Item(157, 11).Value = 1000
$ws.Cells.Item(157, 12).Value = 1000
$ws.Cells.Item(157, 13).Value = 1000
$ws.Cells.Item(157, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(157, 15).Value = "Región del Maule"
$ws.Cells.Item(157, 16).Value = 200
$ws.Cells.Item(157, 17).Value = 5
$ws.Cells.Item(157, 18).Value = "Hortaliza"

# Apply the same date number format as the surrounding rows to the new D157 cell.
$ws.Cells.Item(157, 4).NumberFormat = $ws.Cells.Item(158, 4).NumberFormat

$wb.Save()
